# code thêm tạo report lương tổng hợp
# Update the "last_edited_time" timestamp (column D) for every data row,
# and refresh the updated attendance/work-hour totals exported from Notion
# (Đầy đủ / Tổng công tại CẦN THƠ / Tổng công, etc.) for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data rows share the same last_edited_time string; update them together.
$ws.Range("D2:D22").Value = "2024-08-12T02:00:00.000Z"
$ws.Range("AF2").Value = 11
$ws.Range("AI2").Value = 11
$ws.Range("AM2").Value = 11
$ws.Range("AF3").Value = 11
$ws.Range("AI3").Value = 11
$ws.Range("AM3").Value = 11
$ws.Range("AF4").Value = 11
$ws.Range("AI4").Value = 11
$ws.Range("AM4").Value = 11
$ws.Range("AF5").Value = 11
$ws.Range("AI5").Value = 11
$ws.Range("AM5").Value = 11
$ws.Range("AF8").Value = 11
$ws.Range("AI8").Value = 11
$ws.Range("AM8").Value = 11
$ws.Range("AF9").Value = 11
$ws.Range("AI9").Value = 11
$ws.Range("AM9").Value = 11
$ws.Range("S10").Value = 8
$ws.Range("Y10").Value = 1
$ws.Range("AF10").Value = 7
$ws.Range("AM10").Value = 8
$ws.Range("AF11").Value = 9
$ws.Range("AI11").Value = 9.5
$ws.Range("AM11").Value = 9.5
$ws.Range("AF13").Value = 14
$ws.Range("AI13").Value = 10
$ws.Range("AM13").Value = 14
$ws.Range("AP13").Value = 1
$ws.Range("Y14").Value = 1
$ws.Range("AF14").Value = 7
$ws.Range("AI14").Value = 8
$ws.Range("AM14").Value = 8
$ws.Range("AP14").Value = 1
$ws.Range("S16").Value = 10.5
$ws.Range("AF16").Value = 10
$ws.Range("AM16").Value = 10.5
$ws.Range("AF18").Value = 6
$ws.Range("AI18").Value = 7
$ws.Range("AM18").Value = 7
$ws.Range("S19").Value = 9.5
$ws.Range("AF19").Value = 9
$ws.Range("AM19").Value = 9.5
$ws.Range("AF20").Value = 7
$ws.Range("AI20").Value = 7
$ws.Range("AM20").Value = 7
$ws.Range("S21").Value = 11
$ws.Range("AF21").Value = 11
$ws.Range("AM21").Value = 11
$ws.Range("S22").Value = 10.5
$ws.Range("AF22").Value = 10
$ws.Range("AM22").Value = 10.5
